$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure column D (Price) values are written as text, matching the
# original inline-string cell type, so number-like values (e.g. "217.74")
# are not auto-converted into numeric cells by Excel.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range('D2').Value = '26.960.51'
$ws.Range('E2').Value = '  +0.90%  '
$ws.Range('D3').Value = '1.649.10'
$ws.Range('E3').Value = '  +0.80%  '
$ws.Range('E4').Value = '  -0.07%  '
$ws.Range('D5').Value = '217.74'
$ws.Range('E5').Value = '  +0.29%  '
$ws.Range('E7').Value = '  -0.12%  '
$ws.Range('E8').Value = '  +1.50%  '
$ws.Range('E9').Value = '  +0.10%  '
$ws.Range('D10').Value = '19.79'
$ws.Range('E10').Value = '  +3.99%  '
$ws.Range('D11').Value = '0.0846'
$ws.Range('E11').Value = '  +0.33%  '
$ws.Range('D12').Value = '1.880.26'
$ws.Range('E12').Value = '  +0.85%  '
$ws.Range('D13').Value = '1.641.91'
$ws.Range('E13').Value = '  +0.22%  '
$ws.Range('E14').Value = '  +0.27%  '
$ws.Range('D15').Value = '0.528'
$ws.Range('E15').Value = '  +0.82%  '
$ws.Range('D16').Value = '66.69'
$ws.Range('E16').Value = '  +3.56%  '
$ws.Range('D17').Value = '26.980.64'
$ws.Range('E17').Value = '  +1.04%  '
$ws.Range('D18').Value = '0.0₃0731'
$ws.Range('E18').Value = '  +0.74%  '
$ws.Range('D19').Value = '220.14'
$ws.Range('E19').Value = '  +4.53%  '
$ws.Range('E20').Value = '  -0.04%  '
$ws.Range('E21').Value = '  +1.82%  '
$ws.Range('D22').Value = '6.66'
$ws.Range('E22').Value = '  +7.94%  '
$ws.Range('E23').Value = '  +2.77%  '
$ws.Range('D24').Value = '9.18'
$ws.Range('E24').Value = '  -0.55%  '
$ws.Range('D25').Value = '146.31'
$ws.Range('E25').Value = '  +0.42%  '
$ws.Range('E26').Value = '  -0.05%  '
$ws.Range('E27').Value = '  +4.79%  '
$ws.Range('E28').Value = '  +1.14%  '
$ws.Range('D29').Value = '15.92'
$ws.Range('E29').Value = '  +2.32%  '
$ws.Range('E30').Value = '  +1.55%  '
$ws.Range('E31').Value = '  +0.57%  '
$ws.Range('D32').Value = '3.39'
$ws.Range('E32').Value = '  +1.34%  '
$ws.Range('E33').Value = '  +0.55%  '
$ws.Range('E34').Value = '  +2.39%  '
$ws.Range('D35').Value = '2.46'
$ws.Range('E35').Value = '  +1.06%  '
$ws.Range('D36').Value = '1.250.95'
$ws.Range('E36').Value = '  -1.81%  '
$ws.Range('E37').Value = '  +0.72%  '
$ws.Range('E38').Value = '  +0.63%  '
$ws.Range('D39').Value = '0.831'
$ws.Range('E39').Value = '  +2.88%  '
$ws.Range('E40').Value = '  -0.10%  '
$ws.Range('D41').Value = '0.813'
$ws.Range('E41').Value = '  +1.42%  '
$ws.Range('D42').Value = '5.36'
$ws.Range('E42').Value = '  +2.03%  '
$ws.Range('D43').Value = '1.793.87'
$ws.Range('E43').Value = '  +1.13%  '
$ws.Range('E44').Value = '  -4.83%  '
$ws.Range('D45').Value = '61.33'
$ws.Range('E45').Value = '  +1.35%  '
$ws.Range('D46').Value = '91.55'
$ws.Range('E46').Value = '  +0.56%  '
$ws.Range('E47').Value = '  +1.64%  '
$ws.Range('B48').Value = 'BabyDogeCoin'
$ws.Range('C48').Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range('D48').Value = '0.0₆0105'
$ws.Range('E48').Value = '  +0.36%  '
$ws.Range('B49').Value = 'Cronos'
$ws.Range('C49').Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range('D49').Value = '0.0515'
$ws.Range('E49').Value = '  -0.96%  '
$ws.Range('E50').Value = '  +1.71%  '
$ws.Range('E51').Value = '  +1.21%  '

# Restore the default (Normal) style so column D cells have no explicit
# style index, matching the original workbook formatting.
$ws.Range("D2:D51").Style = "Normal"
